$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A56").Value = "04 Jan 2019"
$ws.Range("B56").Value = 199.9
